$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work")

# --- Row 155-158: unit-length check of the new camera-up / normal vector ---
$ws.Range("A155").Value = "Normal"

$ws.Range("A156").Value = "X"
$ws.Range("B156").Value = -0.70710676900000002
$ws.Range("C156").Formula = "=POWER(B156,2)"

$ws.Range("A157").Value = "Y"
$ws.Range("B157").Value = -0.70710676900000002

$ws.Range("A158").Value = "Z"
$ws.Range("B158").Value = 0

# C157:C158 share one relative formula (C157 is the shared-formula anchor)
$ws.Range("C157:C158").Formula = "=POWER(B157,2)"

$ws.Range("C155").Formula = "=SUM(C156:C158)"

# --- Row 160-161: "Cam Y is:" / "Cam Y Should be:" dot-product check table ---
$ws.Range("A160").Value = "Cam Y is:"
$ws.Range("B160").Value = 3681.8980000000001

$ws.Range("A161").Value = "Cam Y Should be:"
$ws.Range("B161").Value = 923.98900000000003

$ws.Range("D160").Value = "Dot product"
$ws.Range("F160").Value = "mCamUp"
$ws.Range("E160").Value = "tosubject"

# --- Row 159: legend header labels for the H:M helper columns below ---
$ws.Range("H159").Value = "{-3254.375, -1952.625, 609.600}"
$ws.Range("K159").Value = "{-0.707, -0.707, 0.000}"

$ws.Range("H160").Value = -3254.375
$ws.Range("I160").Value = -1952.625
$ws.Range("J160").Value = 609.6
$ws.Range("K160").Value = -0.70699999999999996
$ws.Range("L160").Value = -0.70699999999999996
$ws.Range("M160").Value = 0
$ws.Range("O160").Formula = "=(H160*K160)+(I160*L160)+(J160*M160)"

$ws.Range("H161").Value = -3254.375
$ws.Range("I161").Value = -1952.625
$ws.Range("J161").Value = 609.6
$ws.Range("K161").Value = -0.17499999999999999
$ws.Range("L161").Value = -0.17499999999999999
$ws.Range("M161").Value = 0
$ws.Range("O161").Formula = "=(H161*K161)+(I161*L161)+(J161*M161)"

# --- Scroll the view down to the newly-added rows and leave the next empty cell selected ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 150
$win.ScrollColumn = 1
$ws.Range("K162").Select()

$wb.Save()
